$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsGeneral = $wb.Worksheets.Item(1)
$wsCompany = $wb.Worksheets.Item(2)
$wsUpper   = $wb.Worksheets.Item(3)

$wsGeneral.Name = "duckduck general"
$wsCompany.Name = "duckduck company"
$wsUpper.Name   = "Upper"
